$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "72.329.96"
$ws.Range("E2").Value = "  +3.79%  "
$ws.Range("D3").Value = "2.620.09"
$ws.Range("E3").Value = "  +3.50%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.26"
$ws.Range("E5").Value = "  +0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.47"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.525"
$ws.Range("E8").Value = "  +1.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.173"
$ws.Range("E9").Value = "  +8.90%  "
$ws.Range("D10").Value = "2.619.49"
$ws.Range("E10").Value = "  +3.51%  "
$ws.Range("E11").Value = "  +1.00%  "
$ws.Range("E12").Value = "  +2.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.03"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.104.75"
$ws.Range("E14").Value = "  +4.10%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000188"
$ws.Range("E15").Value = "  +2.86%  "
$ws.Range("D16").Value = "72.307.56"
$ws.Range("E16").Value = "  +3.99%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.64"
$ws.Range("E17").Value = "  +2.08%  "
$ws.Range("D18").Value = "2.622.42"
$ws.Range("E18").Value = "  +5.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.00"
$ws.Range("E19").Value = "  +4.34%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "380.05"
$ws.Range("E20").Value = "  +4.23%  "
$ws.Range("E21").Value = "  +4.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.18"
$ws.Range("E22").Value = "  +2.21%  "
$ws.Range("E23").Value = "  +18.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.32"
$ws.Range("E24").Value = "  +3.81%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("E26").Value = "  +3.16%  "
$ws.Range("E27").Value = "  +9.11%  "
$ws.Range("D28").Value = "2.755.63"
$ws.Range("E28").Value = "  +3.69%  "
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("D30").Value = "0.0₃0952"
$ws.Range("E30").Value = "  +3.09%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.10"
$ws.Range("E31").Value = "  +4.48%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "519.23"
$ws.Range("E32").Value = "  +1.78%  "
$ws.Range("E33").Value = "  +5.95%  "
$ws.Range("E34").Value = "  +2.35%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "165.52"
$ws.Range("E36").Value = "  +2.46%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.31"
$ws.Range("E37").Value = "  +2.72%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.40"
$ws.Range("E38").Value = "  +5.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.07"
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("E40").Value = "  -7.71%  "
$ws.Range("E41").Value = "  +6.08%  "
$ws.Range("E42").Value = "  +4.72%  "
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("E44").Value = "  +6.55%  "
$ws.Range("E45").Value = "  +2.79%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.53"
$ws.Range("E46").Value = "  +1.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "149.67"
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.70"
$ws.Range("E48").Value = "  +2.84%  "
$ws.Range("E49").Value = "  +4.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.69"
$ws.Range("E50").Value = "  +6.41%  "
$ws.Range("D51").Value = "0.0₆0263"
$ws.Range("E51").Value = "  +4.33%  "
